$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate row 4 (the "Test Automation" StorageGrid run) into row 5 so we
#    get an identical second row with the same formatting/styles.
$ws.Range("A4:V4").Copy($ws.Range("A5:V5"))

# 2) Hook up the mailto hyperlink on the new row's username cell, same as A4
#    (adding the hyperlink first - its automatic "Hyperlink" style gets
#    overwritten a couple of lines down when the rest of row 4 is re-applied).
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:APPTESTMBOB4@netapp.com") | Out-Null
$ws.Range("A4:V4").Copy($ws.Range("A5:V5"))

# 3) This new row is a second automation run -> give it its own Opportunity
#    Name value so it reads "Test Second Automation"
$ws.Range("D5").Value = "Test Second Automation"

# 4) Fix the misspelled header "Opprotunity Name" -> "Opportunity Name"
$ws.Range("D3").Value = "Opportunity Name"

# 5) Normalize V3's cell formatting to the same plain header style as the
#    rest of row 3 (drop the redundant fill/border flags it had before).
$ws.Range("U3").Copy()
$ws.Range("V3").PasteSpecial(-4122)

# 6) Widen column D now that it holds the longer "Opportunity Name" /
#    "Test Second Automation" text.
$ws.Columns("D").ColumnWidth = 20.5
